$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Add the two new investment-related parameters that used to be hard-coded
# (maximum investment capacity per year, and max permit build time) to the
# Coupling Parameters table.
$ws.Range("A10").Value = "maximum_investment_capacity_per_year"
$ws.Range("B10").Value = 10000
$ws.Range("A11").Value = "max_permit_build_time"
$ws.Range("B11").Value = 7

# Widen column A so the longer parameter names are fully visible.
$ws.Columns.Item(1).ColumnWidth = 39.83

# Leave a note for the reader about where the authoritative values for the
# new parameter live.
$null = $ws.Range("A11").AddCommentThreaded("check the emlab parameters candidate technologies tab")

# Restore the author's last selection on this sheet.
$ws.Activate()
$null = $ws.Range("H9").Select()
